$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '38.176.19'
$ws.Range("E2").Value = '  +2.60%  '
$ws.Range("D3").Value = '2.056.57'
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").Value = "'229.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.29%  '
$ws.Range("E6").Value = '  +1.09%  '
$ws.Range("D7").Value = "'60.85"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +8.38%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("E9").Value = '  +2.02%  '
$ws.Range("D10").Value = "'0.0805"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.68%  '
$ws.Range("D11").Value = "'0.104"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.03%  '
$ws.Range("D12").Value = "'14.80"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.36%  '
$ws.Range("D13").Value = '2.359.46'
$ws.Range("E13").Value = '  +1.47%  '
$ws.Range("D14").Value = "'21.17"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.10%  '
$ws.Range("E15").Value = '  +3.24%  '
$ws.Range("E16").Value = '  +2.59%  '
$ws.Range("D17").Value = '2.059.59'
$ws.Range("E17").Value = '  +1.81%  '
$ws.Range("D18").Value = '38.108.82'
$ws.Range("E18").Value = '  +2.54%  '
$ws.Range("E19").Value = '  +2.56%  '
$ws.Range("D20").Value = "'69.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.18%  '
$ws.Range("E21").Value = '  +1.45%  '
$ws.Range("D22").Value = "'225.68"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.11%  '
$ws.Range("E23").Value = '  +0.02%  '
$ws.Range("E24").Value = '  -0.10%  '
$ws.Range("D25").Value = "'2.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.44%  '
$ws.Range("D26").Value = "'165.66"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.35%  '
$ws.Range("D27").Value = "'9.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.21%  '
$ws.Range("D28").Value = "'0.134"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.18%  '
$ws.Range("D29").Value = "'18.94"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.99%  '
$ws.Range("E30").Value = '  -0.86%  '
$ws.Range("E31").Value = '  +2.12%  '
$ws.Range("D32").Value = "'4.52"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.35%  '
$ws.Range("E33").Value = '  +2.86%  '
$ws.Range("D34").Value = "'2.05"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +5.31%  '
$ws.Range("E35").Value = '  +0.66%  '
$ws.Range("D36").Value = "'6.46"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +17.58%  '
$ws.Range("E37").Value = '  -2.88%  '
$ws.Range("E38").Value = '  +2.21%  '
$ws.Range("D39").Value = "'1.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.01%  '
$ws.Range("D40").Value = '1.519.33'
$ws.Range("E40").Value = '  +3.13%  '
$ws.Range("B41").Value = 'InjectiveProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D41").Value = "'17.05"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.94%  '
$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").Value = "'97.60"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.58%  '
$ws.Range("E43").Value = '  +0.78%  '
$ws.Range("E44").Value = '  +2.91%  '
$ws.Range("D45").Value = "'0.0924"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.06%  '
$ws.Range("E46").Value = '  +1.66%  '
$ws.Range("D47").Value = "'4.04"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.59%  '
$ws.Range("E48").Value = '  +0.59%  '
$ws.Range("E49").Value = '  +1.45%  '
$ws.Range("D50").Value = "'7.04"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.80%  '
$ws.Range("D51").Value = '2.248.42'
$ws.Range("E51").Value = '  +1.60%  '
